$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$nm = $p.NotesMaster
Write-Host "SlideMaster CustomLayouts Count:"
Write-Host $m.CustomLayouts.Count
Write-Host "NotesMaster CustomLayouts Count:"
Write-Host $nm.CustomLayouts.Count
